$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 143
$ws.Range("B143").Value = 6937247
$ws.Range("E143").Value = "AEK Athens"
$ws.Range("F143").Value = "Asteras Tripolis"
$ws.Range("G143").Value = 4
$ws.Range("H143").Value = 2
$ws.Range("I143").Value = "H"
$ws.Range("J143").Value = 1.285
$ws.Range("K143").Value = 5.5
$ws.Range("L143").Value = 12
$ws.Range("M143").Value = 1.285
$ws.Range("N143").Value = 5.75
$ws.Range("O143").Value = 10
$ws.Range("P143").Value = -1.5
$ws.Range("Q143").Value = 1.825
$ws.Range("R143").Value = 2.025
$ws.Range("S143").Value = 3
$ws.Range("T143").Value = 2.025
$ws.Range("U143").Value = 1.825
$ws.Range("V143").Value = 0.2849999999999999
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 0.825
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 1.025

# Row 144
$ws.Range("B144").Value = 6937250
$ws.Range("E144").Value = "Giannina"
$ws.Range("F144").Value = "Lamia"
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 4
$ws.Range("I144").Value = "A"
$ws.Range("J144").Value = 2.3
$ws.Range("K144").Value = 3.25
$ws.Range("L144").Value = 3.25
$ws.Range("M144").Value = 2.55
$ws.Range("N144").Value = 2.875
$ws.Range("O144").Value = 3.1
$ws.Range("P144").Value = 0
$ws.Range("Q144").Value = 1.75
$ws.Range("R144").Value = 2.125
$ws.Range("S144").Value = 2
$ws.Range("T144").Value = 1.85
$ws.Range("U144").Value = 2
$ws.Range("V144").Value = -1
$ws.Range("X144").Value = 2.1
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = 1.125
$ws.Range("AA144").Value = 0.8500000000000001

# Row 170
$ws.Range("B170").Value = 6937266
$ws.Range("E170").Value = "Atromitos Athinon"
$ws.Range("F170").Value = "Lamia"
$ws.Range("G170").Value = 3
$ws.Range("H170").Value = 1
$ws.Range("I170").Value = "H"
$ws.Range("J170").Value = 2.3
$ws.Range("K170").Value = 3.2
$ws.Range("L170").Value = 3.1
$ws.Range("M170").Value = 2.2
$ws.Range("N170").Value = 3.3
$ws.Range("O170").Value = 3.3
$ws.Range("P170").Value = -0.25
$ws.Range("Q170").Value = 1.925
$ws.Range("R170").Value = 1.925
$ws.Range("S170").Value = 2.5
$ws.Range("T170").Value = 2.025
$ws.Range("U170").Value = 1.825
$ws.Range("V170").Value = 1.2
$ws.Range("X170").Value = -1
$ws.Range("Y170").Value = 0.925
$ws.Range("Z170").Value = -1
$ws.Range("AA170").Value = 1.025
$ws.Range("AB170").Value = -1

# Row 171
$ws.Range("B171").Value = 6937268
$ws.Range("E171").Value = "Panetolikos"
$ws.Range("F171").Value = "Olympiakos"
$ws.Range("G171").Value = 1
$ws.Range("H171").Value = 2
$ws.Range("I171").Value = "A"
$ws.Range("J171").Value = 8
$ws.Range("K171").Value = 5
$ws.Range("L171").Value = 1.363
$ws.Range("M171").Value = 8.5
$ws.Range("N171").Value = 5
$ws.Range("O171").Value = 1.363
$ws.Range("P171").Value = 1.25
$ws.Range("Q171").Value = 2.025
$ws.Range("R171").Value = 1.825
$ws.Range("S171").Value = 2.75
$ws.Range("T171").Value = 1.85
$ws.Range("U171").Value = 2
$ws.Range("V171").Value = -1
$ws.Range("X171").Value = 0.363
$ws.Range("Y171").Value = 0.5125
$ws.Range("Z171").Value = -0.5
$ws.Range("AA171").Value = 0.425
$ws.Range("AB171").Value = -0.5

# Row 177
$ws.Range("B177").Value = 6937270
$ws.Range("E177").Value = "Olympiakos"
$ws.Range("F177").Value = "Volos NFC"
$ws.Range("G177").Value = 3
$ws.Range("H177").Value = 0
$ws.Range("J177").Value = 1.125
$ws.Range("K177").Value = 9
$ws.Range("L177").Value = 19
$ws.Range("M177").Value = 1.111
$ws.Range("N177").Value = 9
$ws.Range("O177").Value = 21
$ws.Range("P177").Value = -2.25
$ws.Range("Q177").Value = 1.875
$ws.Range("R177").Value = 1.975
$ws.Range("S177").Value = 3.25
$ws.Range("V177").Value = 0.111
$ws.Range("Y177").Value = 0.875
$ws.Range("AA177").Value = -0.5
$ws.Range("AB177").Value = 0.425

# Row 178
$ws.Range("B178").Value = 6937272
$ws.Range("E178").Value = "Lamia"
$ws.Range("F178").Value = "PAOK Salonika"
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 2
$ws.Range("I178").Value = "A"
$ws.Range("J178").Value = 7.5
$ws.Range("K178").Value = 4.5
$ws.Range("L178").Value = 1.444
$ws.Range("M178").Value = 9.5
$ws.Range("N178").Value = 5
$ws.Range("O178").Value = 1.333
$ws.Range("P178").Value = 1.5
$ws.Range("Q178").Value = 1.925
$ws.Range("R178").Value = 1.925
$ws.Range("S178").Value = 3
$ws.Range("T178").Value = 1.95
$ws.Range("U178").Value = 1.9
$ws.Range("V178").Value = -1
$ws.Range("X178").Value = 0.333
$ws.Range("Y178").Value = -1
$ws.Range("Z178").Value = 0.925
$ws.Range("AA178").Value = -1
$ws.Range("AB178").Value = 0.8999999999999999

# Row 179
$ws.Range("B179").Value = 6935700
$ws.Range("E179").Value = "Panserraikos"
$ws.Range("F179").Value = "Asteras Tripolis"
$ws.Range("G179").Value = 2
$ws.Range("H179").Value = 1
$ws.Range("I179").Value = "H"
$ws.Range("J179").Value = 2.6
$ws.Range("K179").Value = 3.2
$ws.Range("L179").Value = 2.875
$ws.Range("M179").Value = 2.25
$ws.Range("N179").Value = 3.3
$ws.Range("O179").Value = 3.3
$ws.Range("P179").Value = -0.25
$ws.Range("S179").Value = 2.25
$ws.Range("T179").Value = 2
$ws.Range("U179").Value = 1.85
$ws.Range("V179").Value = 1.25
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = 0.925
$ws.Range("Z179").Value = -1
$ws.Range("AA179").Value = 1
$ws.Range("AB179").Value = -1

# Row 194
$ws.Range("B194").Value = 7920470
$ws.Range("E194").Value = "AEK Athens"
$ws.Range("F194").Value = "Olympiakos"
$ws.Range("G194").Value = 1
$ws.Range("H194").Value = 0
$ws.Range("J194").Value = 1.909
$ws.Range("K194").Value = 3.4
$ws.Range("L194").Value = 4.2
$ws.Range("M194").Value = 2.2
$ws.Range("N194").Value = 3.2
$ws.Range("O194").Value = 3.5
$ws.Range("P194").Value = -0.25
$ws.Range("Q194").Value = 1.85
$ws.Range("R194").Value = 2
$ws.Range("S194").Value = 2.5
$ws.Range("V194").Value = 1.2
$ws.Range("Y194").Value = 0.8500000000000001
$ws.Range("AA194").Value = -1
$ws.Range("AB194").Value = 0.825

# Row 195
$ws.Range("B195").Value = 7920471
$ws.Range("E195").Value = "Aris Salonika"
$ws.Range("F195").Value = "Lamia"
$ws.Range("G195").Value = 3
$ws.Range("H195").Value = 1
$ws.Range("J195").Value = 1.571
$ws.Range("K195").Value = 4
$ws.Range("L195").Value = 6
$ws.Range("M195").Value = 1.444
$ws.Range("N195").Value = 4.5
$ws.Range("O195").Value = 8.5
$ws.Range("P195").Value = -1.25
$ws.Range("Q195").Value = 1.925
$ws.Range("R195").Value = 1.925
$ws.Range("S195").Value = 2.75
$ws.Range("V195").Value = 0.444
$ws.Range("Y195").Value = 0.925
$ws.Range("AA195").Value = 1.025
$ws.Range("AB195").Value = -1

# Row 237
$ws.Range("B237").NumberFormat = "@"
$ws.Range("B237").Value = "8140565"
$ws.Range("B237").Style = "Normal"
$ws.Range("E237").Value = "Panathinaikos"
$ws.Range("F237").Value = "Olympiakos"
$ws.Range("J237").Value = 2.4
$ws.Range("K237").Value = 3.3
$ws.Range("L237").Value = 2.8
$ws.Range("M237").Value = 2.75
$ws.Range("N237").Value = 3.3
$ws.Range("O237").Value = 2.5
$ws.Range("P237").Value = 0
$ws.Range("Q237").Value = 2.025
$ws.Range("R237").Value = 1.825
$ws.Range("S237").Value = 2.5
$ws.Range("T237").Value = 1.85
$ws.Range("U237").Value = 2

# Row 239
$ws.Range("B239").NumberFormat = "@"
$ws.Range("B239").Value = "8140226"
$ws.Range("B239").Style = "Normal"
$ws.Range("E239").Value = "Aris Salonika"
$ws.Range("F239").Value = "PAOK Salonika"
$ws.Range("J239").Value = 5.75
$ws.Range("K239").Value = 4.5
$ws.Range("L239").Value = 1.5
$ws.Range("M239").Value = 4.5
$ws.Range("N239").Value = 4.1
$ws.Range("O239").Value = 1.666
$ws.Range("P239").Value = 0.75
$ws.Range("Q239").Value = 1.975
$ws.Range("R239").Value = 1.875
$ws.Range("S239").Value = 2.75
$ws.Range("T239").Value = 2
$ws.Range("U239").Value = 1.85
